$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s16 = @'
Criar tabela SegmentoClassificacao
'@

$s17 = @'
Preciso gerar um script para criar uma tabela no SQL Server onde deverá constar a primary key identity ID, a sigla e o descritivo.
'@

$s18 = @'
Criar tabela Segmento
'@

$s19 = @'
Preciso gerar um script para criar uma tabela no SQL Server com o nome "Setor Econômico" onde deverá constar a primary key identity ID e o descritivo.
'@

$s20 = @'
Criar tabela Subsetor
'@

$s21 = @'
Preciso gerar um script para criar uma tabela no SQL Server com o nome "Subsetor" onde deverá constar a primary key identity ID e o descritivo.
'@

$s22 = @'
Criar tabela SetorEconomico
'@

$s23 = @'
Preciso gerar um script para criar uma tabela no SQL Server com o nome "Segmento" onde deverá constar a primary key identity ID e o descritivo.
'@

$s24 = @'
Criar tabela Empresa
'@

$s25 = @'
Preciso gerar um script para criar uma tabela no SQL Server com o nome "Empresa" onde deverá constar a primary key identity ID, Nome, código (sigla) fazendo chave estrangeira para as seguintes 
tabelas dbo.SegmentoClassificacao (NÃO obrigatória), SetorEconomico (obrigatória), Subsetor (obrigatória) e Segmento (obrigatória).
'@

$s26 = @'
Para criar uma tabela no SQL Server com uma **primary key identity ID**, uma coluna para a **sigla** e outra para o **descritivo**, você pode usar um script SQL. Abaixo está um exemplo de como criar essa tabela:
---
### Script SQL para Criar a Tabela
```sql
CREATE TABLE Ativos (
    ID INT IDENTITY(1,1) PRIMARY KEY,  -- Coluna ID como primary key e autoincremento
    Sigla NVARCHAR(10) NOT NULL,       -- Coluna para a sigla (tamanho máximo de 10 caracteres)
    Descritivo NVARCHAR(100) NOT NULL  -- Coluna para o descritivo (tamanho máximo de 100 caracteres)
);
```
---
### Explicação do Script
1. **`ID INT IDENTITY(1,1) PRIMARY KEY`**:
   - `ID`: Nome da coluna.
   - `INT`: Tipo de dado inteiro.
   - `IDENTITY(1,1)`: Define que a coluna é autoincrementada, começando em 1 e incrementando de 1 em 1.
   - `PRIMARY KEY`: Define a coluna como chave primária.
2. **`Sigla NVARCHAR(10) NOT NULL`**:
   - `Sigla`: Nome da coluna.
   - `NVARCHAR(10)`: Tipo de dado para strings com tamanho máximo de 10 caracteres.
   - `NOT NULL`: Impede que o valor seja nulo.
3. **`Descritivo NVARCHAR(100) NOT NULL`**:
   - `Descritivo`: Nome da coluna.
   - `NVARCHAR(100)`: Tipo de dado para strings com tamanho máximo de 100 caracteres.
   - `NOT NULL`: Impede que o valor seja nulo.
---
### Como Executar o Script
1. **No SQL Server Management Studio (SSMS)**:
   - Abra o SSMS.
   - Conecte-se ao banco de dados desejado.
   - Abra uma nova janela de consulta.
   - Cole o script SQL acima.
   - Execute o script (tecla `F5` ou botão "Executar").
2. **Via Python (usando SQLAlchemy ou pyodbc)**:
   Se você quiser executar o script diretamente do Python, pode usar uma biblioteca como `pyodbc` ou `SQLAlchemy`. Aqui está um exemplo usando `pyodbc`:
   ```python
   import pyodbc
   # Configurações de conexão
   server = "localhost"
   database = "SeuBancoDeDados"
   username = "sa"
   password = "YourPassword123"
   connection_string = f"DRIVER={{ODBC Driver 17 for SQL Server}};SERVER={server};DATABASE={database};UID={username};PWD={password}"
   # Conecta ao banco de dados
   conn = pyodbc.connect(connection_string)
   cursor = conn.cursor()
   # Script SQL para criar a tabela
   create_table_sql = """
   CREATE TABLE Ativos (
       ID INT IDENTITY(1,1) PRIMARY KEY,
       Sigla NVARCHAR(10) NOT NULL,
       Descritivo NVARCHAR(100) NOT NULL
   );
   """
   # Executa o script
   cursor.execute(create_table_sql)
   conn.commit()
   # Fecha a conexão
   cursor.close()
   conn.close()
   print("Tabela 'Ativos' criada com sucesso!")
   ```
3. **Via Docker (se estiver usando um container SQL Server)**:
   Se você estiver usando um container Docker para o SQL Server, pode executar o script diretamente no container:
   ```bash
   docker exec -it sqlserver_db /opt/mssql-tools/bin/sqlcmd -S localhost -U sa -P YourPassword123 -d SeuBancoDeDados -Q "CREATE TABLE Ativos (ID INT IDENTITY(1,1) PRIMARY KEY, Sigla NVARCHAR(10) NOT NULL, Descritivo NVARCHAR(100) NOT NULL);"
   ```
---
### Verificando a Tabela Criada
Após executar o script, você pode verificar se a tabela foi criada com sucesso:
1. No SSMS:
   - Expanda o banco de dados no Object Explorer.
   - Verifique se a tabela `Ativos` aparece na lista de tabelas.
2. Via SQL:
   - Execute a consulta:
     ```sql
     SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Ativos';
     ```
   - Se a tabela existir, ela será listada.
3. Via Python:
   - Use o seguinte código para listar as tabelas:
     ```python
     cursor.execute("SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Ativos';")
     rows = cursor.fetchall()
     for row in rows:
         print(row)
     ```
---
### Adicionando Dados à Tabela
Agora que a tabela está criada, você pode inserir dados nela. Aqui está um exemplo de inserção:
```sql
INSERT INTO Ativos (Sigla, Descritivo)
VALUES ('PETR4', 'Petrobras PN'),
       ('VALE3', 'Vale ON'),
       ('ITUB4', 'Itaú Unibanco PN');
```
Ou, via Python:
```python
insert_sql = """
INSERT INTO Ativos (Sigla, Descritivo)
VALUES (?, ?);
"""
data = [
    ('PETR4', 'Petrobras PN'),
    ('VALE3', 'Vale ON'),
    ('ITUB4', 'Itaú Unibanco PN'),
]
cursor.executemany(insert_sql, data)
conn.commit()
```
---
### Resumo
- O script SQL cria uma tabela `Ativos` com uma coluna `ID` autoincrementada como chave primária, uma coluna `Sigla` e uma coluna `Descritivo`.
- Você pode executar o script diretamente no SQL Server, via Python ou em um container Docker.
- Após criar a tabela, você pode inserir dados e consultá-los.
Com isso, sua tabela estará pronta para uso no seu projeto!
'@

$s27 = @'
Para criar uma tabela no SQL Server chamada **"Setor Econômico"** com uma **primary key identity ID** e uma coluna para o **descritivo**, você pode usar o seguinte script SQL:
---
### Script SQL para Criar a Tabela
```sql
CREATE TABLE [Setor Econômico] (
    ID INT IDENTITY(1,1) PRIMARY KEY,  -- Coluna ID como primary key e autoincremento
    Descritivo NVARCHAR(100) NOT NULL  -- Coluna para o descritivo (tamanho máximo de 100 caracteres)
);
```
---
### Explicação do Script
1. **`ID INT IDENTITY(1,1) PRIMARY KEY`**:
   - `ID`: Nome da coluna.
   - `INT`: Tipo de dado inteiro.
   - `IDENTITY(1,1)`: Define que a coluna é autoincrementada, começando em 1 e incrementando de 1 em 1.
   - `PRIMARY KEY`: Define a coluna como chave primária.
2. **`Descritivo NVARCHAR(100) NOT NULL`**:
   - `Descritivo`: Nome da coluna.
   - `NVARCHAR(100)`: Tipo de dado para strings com tamanho máximo de 100 caracteres.
   - `NOT NULL`: Impede que o valor seja nulo.
3. **`[Setor Econômico]`**:
   - O nome da tabela está entre colchetes (`[ ]`) porque contém espaços. No SQL Server, nomes de tabelas ou colunas com espaços ou caracteres especiais devem ser escapados com colchetes.
---
### Como Executar o Script
1. **No SQL Server Management Studio (SSMS)**:
   - Abra o SSMS.
   - Conecte-se ao banco de dados desejado.
   - Abra uma nova janela de consulta.
   - Cole o script SQL acima.
   - Execute o script (tecla `F5` ou botão "Executar").
2. **Via Python (usando SQLAlchemy ou pyodbc)**:
   Se você quiser executar o script diretamente do Python, pode usar uma biblioteca como `pyodbc` ou `SQLAlchemy`. Aqui está um exemplo usando `pyodbc`:
   ```python
   import pyodbc
   # Configurações de conexão
   server = "localhost"
   database = "SeuBancoDeDados"
   username = "sa"
   password = "YourPassword123"
   connection_string = f"DRIVER={{ODBC Driver 17 for SQL Server}};SERVER={server};DATABASE={database};UID={username};PWD={password}"
   # Conecta ao banco de dados
   conn = pyodbc.connect(connection_string)
   cursor = conn.cursor()
   # Script SQL para criar a tabela
   create_table_sql = """
   CREATE TABLE [Setor Econômico] (
       ID INT IDENTITY(1,1) PRIMARY KEY,
       Descritivo NVARCHAR(100) NOT NULL
   );
   """
   # Executa o script
   cursor.execute(create_table_sql)
   conn.commit()
   # Fecha a conexão
   cursor.close()
   conn.close()
   print("Tabela 'Setor Econômico' criada com sucesso!")
   ```
3. **Via Docker (se estiver usando um container SQL Server)**:
   Se você estiver usando um container Docker para o SQL Server, pode executar o script diretamente no container:
   ```bash
   docker exec -it sqlserver_db /opt/mssql-tools/bin/sqlcmd -S localhost -U sa -P YourPassword123 -d SeuBancoDeDados -Q "CREATE TABLE [Setor Econômico] (ID INT IDENTITY(1,1) PRIMARY KEY, Descritivo NVARCHAR(100) NOT NULL);"
   ```
---
### Verificando a Tabela Criada
Após executar o script, você pode verificar se a tabela foi criada com sucesso:
1. No SSMS:
   - Expanda o banco de dados no Object Explorer.
   - Verifique se a tabela `Setor Econômico` aparece na lista de tabelas.
2. Via SQL:
   - Execute a consulta:
     ```sql
     SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Setor Econômico';
     ```
   - Se a tabela existir, ela será listada.
3. Via Python:
   - Use o seguinte código para listar as tabelas:
     ```python
     cursor.execute("SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Setor Econômico';")
     rows = cursor.fetchall()
     for row in rows:
         print(row)
     ```
---
### Adicionando Dados à Tabela
Agora que a tabela está criada, você pode inserir dados nela. Aqui está um exemplo de inserção:
```sql
INSERT INTO [Setor Econômico] (Descritivo)
VALUES ('Tecnologia da Informação'),
       ('Energia'),
       ('Financeiro');
```
Ou, via Python:
```python
insert_sql = """
INSERT INTO [Setor Econômico] (Descritivo)
VALUES (?);
"""
data = [
    ('Tecnologia da Informação',),
    ('Energia',),
    ('Financeiro',),
]
cursor.executemany(insert_sql, data)
conn.commit()
```
---
### Resumo
- O script SQL cria uma tabela chamada `Setor Econômico` com uma coluna `ID` autoincrementada como chave primária e uma coluna `Descritivo`.
- Você pode executar o script diretamente no SQL Server, via Python ou em um container Docker.
- Após criar a tabela, você pode inserir dados e consultá-los.
Com isso, sua tabela `Setor Econômico` estará pronta para uso no seu projeto!
'@

$s28 = @'
Para criar uma tabela no SQL Server chamada **"Subsetor"** com uma **primary key identity ID** e uma coluna para o **descritivo**, você pode usar o seguinte script SQL:
---
### Script SQL para Criar a Tabela
```sql
CREATE TABLE Subsetor (
    ID INT IDENTITY(1,1) PRIMARY KEY,  -- Coluna ID como primary key e autoincremento
    Descritivo NVARCHAR(100) NOT NULL  -- Coluna para o descritivo (tamanho máximo de 100 caracteres)
);
```
---
### Explicação do Script
1. **`ID INT IDENTITY(1,1) PRIMARY KEY`**:
   - `ID`: Nome da coluna.
   - `INT`: Tipo de dado inteiro.
   - `IDENTITY(1,1)`: Define que a coluna é autoincrementada, começando em 1 e incrementando de 1 em 1.
   - `PRIMARY KEY`: Define a coluna como chave primária.
2. **`Descritivo NVARCHAR(100) NOT NULL`**:
   - `Descritivo`: Nome da coluna.
   - `NVARCHAR(100)`: Tipo de dado para strings com tamanho máximo de 100 caracteres.
   - `NOT NULL`: Impede que o valor seja nulo.
3. **`Subsetor`**:
   - Nome da tabela. Como não contém espaços ou caracteres especiais, não é necessário usar colchetes.
---
### Como Executar o Script
1. **No SQL Server Management Studio (SSMS)**:
   - Abra o SSMS.
   - Conecte-se ao banco de dados desejado.
   - Abra uma nova janela de consulta.
   - Cole o script SQL acima.
   - Execute o script (tecla `F5` ou botão "Executar").
2. **Via Python (usando SQLAlchemy ou pyodbc)**:
   Se você quiser executar o script diretamente do Python, pode usar uma biblioteca como `pyodbc` ou `SQLAlchemy`. Aqui está um exemplo usando `pyodbc`:
   ```python
   import pyodbc
   # Configurações de conexão
   server = "localhost"
   database = "SeuBancoDeDados"
   username = "sa"
   password = "YourPassword123"
   connection_string = f"DRIVER={{ODBC Driver 17 for SQL Server}};SERVER={server};DATABASE={database};UID={username};PWD={password}"
   # Conecta ao banco de dados
   conn = pyodbc.connect(connection_string)
   cursor = conn.cursor()
   # Script SQL para criar a tabela
   create_table_sql = """
   CREATE TABLE Subsetor (
       ID INT IDENTITY(1,1) PRIMARY KEY,
       Descritivo NVARCHAR(100) NOT NULL
   );
   """
   # Executa o script
   cursor.execute(create_table_sql)
   conn.commit()
   # Fecha a conexão
   cursor.close()
   conn.close()
   print("Tabela 'Subsetor' criada com sucesso!")
   ```
3. **Via Docker (se estiver usando um container SQL Server)**:
   Se você estiver usando um container Docker para o SQL Server, pode executar o script diretamente no container:
   ```bash
   docker exec -it sqlserver_db /opt/mssql-tools/bin/sqlcmd -S localhost -U sa -P YourPassword123 -d SeuBancoDeDados -Q "CREATE TABLE Subsetor (ID INT IDENTITY(1,1) PRIMARY KEY, Descritivo NVARCHAR(100) NOT NULL);"
   ```
---
### Verificando a Tabela Criada
Após executar o script, você pode verificar se a tabela foi criada com sucesso:
1. No SSMS:
   - Expanda o banco de dados no Object Explorer.
   - Verifique se a tabela `Subsetor` aparece na lista de tabelas.
2. Via SQL:
   - Execute a consulta:
     ```sql
     SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Subsetor';
     ```
   - Se a tabela existir, ela será listada.
3. Via Python:
   - Use o seguinte código para listar as tabelas:
     ```python
     cursor.execute("SELECT * FROM INFORMATION_SCHEMA.TABLES WHERE TABLE_NAME = 'Subsetor';")
     rows = cursor.fetchall()
     for row in rows:
         print(row)
     ```
---
### Adicionando Dados à Tabela
Agora que a tabela está criada, você pode inserir dados nela. Aqui está um exemplo de inserção:
```sql
INSERT INTO Subsetor (Descritivo)
VALUES ('Software'),
       ('Petróleo e Gás'),
       ('Bancos');
```
Ou, via Python:
```python
insert_sql = """
INSERT INTO Subsetor (Descritivo)
VALUES (?);
"""
data = [
    ('Software',),
    ('Petróleo e Gás',),
    ('Bancos',),
]
cursor.executemany(insert_sql, data)
conn.commit()
```
---
### Resumo
- O script SQL cria uma tabela chamada `Subsetor` com uma coluna `ID` autoincrementada como chave primária e uma coluna `Descritivo`.
- Você pode executar o script diretamente no SQL Server, via Python ou em um container Docker.
- Após criar a tabela, você pode inserir dados e consultá-los.
Com isso, sua tabela `Subsetor` estará pronta para uso no seu projeto!
'@

$s29 = @'
O servidor não respondeu a tempo então tive que fazer a query manualmente
'@

# Step 1: Fill in Objetivo (A) and Pergunta (B) columns for all new rows, in row order
$ws.Range("A5").Value = $s16
$ws.Range("B5").Value = $s17
$ws.Range("A6").Value = $s18
$ws.Range("B6").Value = $s19
$ws.Range("A7").Value = $s20
$ws.Range("B7").Value = $s21
$ws.Range("A8").Value = $s22
$ws.Range("B8").Value = $s23
$ws.Range("A9").Value = $s24
$ws.Range("B9").Value = $s25

# Step 2: Fill in Resposta (C) column for rows 5-7
$ws.Range("C5").Value = $s26
$ws.Range("C6").Value = $s27
$ws.Range("C7").Value = $s28

# Step 3: Fill in Observação (D) column for row 8
$ws.Range("D8").Value = $s29

# Step 4: Time values (Tempo início / Tempo fim)
$ws.Range("E5").Value = 0.94305555555555554
$ws.Range("F5").Value = 0.9458333333333333
$ws.Range("E6").Value = 0.9458333333333333
$ws.Range("F6").Value = 0.94722222222222219
$ws.Range("E7").Value = 0.94791666666666663
$ws.Range("F7").Value = 0.95
$ws.Range("E8").Value = 0.96180555555555558
$ws.Range("F8").Value = 0.96388888888888891
$ws.Range("E5:F8").NumberFormat = "h:mm"

# Step 5: Formatting - wrap text for long-form cells
$ws.Range("C5").WrapText = $true
$ws.Range("C6").WrapText = $true
$ws.Range("C7").WrapText = $true
$ws.Range("B9").WrapText = $true

# Step 6: Row heights
$ws.Rows.Item(5).RowHeight = 409.5
$ws.Rows.Item(6).RowHeight = 409.5
$ws.Rows.Item(7).RowHeight = 409.5
$ws.Rows.Item(9).RowHeight = 120

# Step 7: View / selection updates
$excel.ActiveWindow.Zoom = 85
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F8").Select()
